# Update the "Facebook" worksheet (xl/worksheets/sheet1.xml) with two new
# columns (N = ratio of literals/variables, O = variables, P = literals)
# and make the Facebook sheet the active/selected sheet instead of the
# "Chart3" chartsheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Facebook")

# Make "Facebook" the active sheet (this also clears tabSelected on the
# previously active chartsheet "Chart3").
$ws.Activate()

# New header cells for columns O and P - copy the header style (bold /
# shaded fill) from the existing header cell A2, then set the text.
$ws.Range("A2").Copy()
$ws.Range("O2").PasteSpecial(-4122)
$ws.Range("O2").Value = "variables"

$ws.Range("A2").Copy()
$ws.Range("P2").PasteSpecial(-4122)
$ws.Range("P2").Value = "literals"

# Column N: ratio of literals (L) to variables (J) for rows 3-12.
$ws.Range("N3").Formula = "=L3/J3"
$ws.Range("N4:N12").Formula = "=L4/J4"

# Column O: variable counts for rows 3-12.
$ws.Range("O3").Value = 2286
$ws.Range("O4").Value = 7193
$ws.Range("O5").Value = 1559
$ws.Range("O6").Value = 1041
$ws.Range("O7").Value = 1154
$ws.Range("O8").Value = 424
$ws.Range("O9").Value = 5457
$ws.Range("O10").Value = 5220
$ws.Range("O11").Value = 3644
$ws.Range("O12").Value = 355

# Column P: literal counts for rows 3-12.
# P3 additionally carries the red-highlight style used elsewhere in the
# sheet (e.g. row 13 / A13), so copy that formatting across first.
$ws.Range("A13").Copy()
$ws.Range("P3").PasteSpecial(-4122)
$ws.Range("P3").Value = 60199

$ws.Range("P4").Value = 531836
$ws.Range("P5").Value = 26819
$ws.Range("P6").Value = 12909
$ws.Range("P7").Value = 16091
$ws.Range("P8").Value = 2931
$ws.Range("P9").Value = 312068
$ws.Range("P10").Value = 265496
$ws.Range("P11").Value = 149737
$ws.Range("P12").Value = 2287

# Move the active cell/selection on the Facebook sheet from G26 to G27.
[void]$ws.Range("G27").Select()
